# Applies the "New crime data collected" weekly refresh:
#   - bump the report Volume/Number and the covered week dates
#   - update the crime-complaint figures (current/prior period counts
#     and their computed % changes) for rows 14-33
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report week dates) ---
$ws.Range("A8").Value = "Volume 31   Number  12"
$ws.Range("C9").Value = "Report Covering the Week  3/18/2024  Through  3/24/2024"

# --- Crime complaint statistics table updates (rows 14-33) ---
# Row 14
$ws.Range("C14").Value = 2
$ws.Range("F14").Value = 5
$ws.Range("G14").Value = 5
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 10
$ws.Range("K14").Value = -28.571428571428
$ws.Range("L14").Value = -16.666666666666
$ws.Range("M14").Value = -9.090909090909
$ws.Range("N14").Value = -86.842105263157

# Row 15
$ws.Range("C15").Value = 1
$ws.Range("E15").Value = -50
$ws.Range("F15").Value = 6
$ws.Range("G15").Value = 10
$ws.Range("H15").Value = -40
$ws.Range("I15").Value = 27
$ws.Range("J15").Value = 31
$ws.Range("K15").Value = -12.903225806451
$ws.Range("L15").Value = -30.769230769230
$ws.Range("M15").Value = -34.146341463414
$ws.Range("N15").Value = -73

# Row 16
$ws.Range("C16").Value = 33
$ws.Range("D16").Value = 40
$ws.Range("E16").Value = -17.5
$ws.Range("F16").Value = 145
$ws.Range("H16").Value = 16
$ws.Range("I16").Value = 451
$ws.Range("J16").Value = 408
$ws.Range("K16").Value = 10.539215686274
$ws.Range("L16").Value = 2.968036529680
$ws.Range("M16").Value = -11.045364891518
$ws.Range("N16").Value = -78.117418728772

# Row 17
$ws.Range("C17").Value = 54
$ws.Range("D17").Value = 56
$ws.Range("E17").Value = -3.571428571428
$ws.Range("F17").Value = 236
$ws.Range("G17").Value = 208
$ws.Range("H17").Value = 13.461538461538
$ws.Range("I17").Value = 644
$ws.Range("J17").Value = 616
$ws.Range("K17").Value = 4.545454545454
$ws.Range("L17").Value = 9.523809523809
$ws.Range("M17").Value = 64.705882352941
$ws.Range("N17").Value = -40.425531914893

# Row 18
$ws.Range("C18").Value = 22
$ws.Range("D18").Value = 32
$ws.Range("E18").Value = -31.25
$ws.Range("F18").Value = 112
$ws.Range("G18").Value = 139
$ws.Range("H18").Value = -19.424460431654
$ws.Range("I18").Value = 325
$ws.Range("J18").Value = 387
$ws.Range("K18").Value = -16.020671834625
$ws.Range("L18").Value = -25.629290617849
$ws.Range("M18").Value = -1.515151515151
$ws.Range("N18").Value = -86.761710794297

# Row 19
$ws.Range("C19").Value = 125
$ws.Range("D19").Value = 140
$ws.Range("E19").Value = -10.714285714285
$ws.Range("F19").Value = 473
$ws.Range("G19").Value = 450
$ws.Range("H19").Value = 5.111111111111
$ws.Range("I19").Value = 1377
$ws.Range("J19").Value = 1303
$ws.Range("K19").Value = 5.679201841903
$ws.Range("L19").Value = 4.397270659590
$ws.Range("M19").Value = 41.086065573770
$ws.Range("N19").Value = -42.481203007518

# Row 20
$ws.Range("C20").Value = 20
$ws.Range("D20").Value = 18
$ws.Range("E20").Value = 11.111111111111
$ws.Range("F20").Value = 67
$ws.Range("G20").Value = 78
$ws.Range("H20").Value = -14.102564102564
$ws.Range("I20").Value = 203
$ws.Range("J20").Value = 259
$ws.Range("K20").Value = -21.621621621621
$ws.Range("L20").Value = -29.757785467128
$ws.Range("M20").Value = 123.076923076923
$ws.Range("N20").Value = -90.627885503231

# Row 21
$ws.Range("C21").Value = 257
$ws.Range("D21").Value = 288
$ws.Range("E21").Value = -10.763888888888
$ws.Range("F21").Value = 1044
$ws.Range("G21").Value = 1015
$ws.Range("H21").Value = 2.857142857142
$ws.Range("I21").Value = 3037
$ws.Range("J21").Value = 3018
$ws.Range("K21").Value = 0.629555997349
$ws.Range("L21").Value = -2.722613709160
$ws.Range("M21").Value = 29.399233063485
$ws.Range("N21").Value = -70.608729313848

# Row 22
$ws.Range("C22").Value = 5
$ws.Range("D22").Value = 4
$ws.Range("E22").Value = 25
$ws.Range("F22").Value = 16
$ws.Range("G22").Value = 18
$ws.Range("H22").Value = -11.111111111111
$ws.Range("I22").Value = 62
$ws.Range("J22").Value = 63
$ws.Range("K22").Value = -1.587301587301
$ws.Range("L22").Value = -21.518987341772
$ws.Range("M22").Value = 14.814814814814

# Row 23
$ws.Range("C23").Value = 16
$ws.Range("D23").Value = 27
$ws.Range("E23").Value = -40.740740740740
$ws.Range("F23").Value = 94
$ws.Range("G23").Value = 103
$ws.Range("H23").Value = -8.737864077669
$ws.Range("I23").Value = 288
$ws.Range("J23").Value = 271
$ws.Range("K23").Value = 6.273062730627
$ws.Range("L23").Value = 11.196911196911
$ws.Range("M23").Value = 60

# Row 24
$ws.Range("C24").Value = 247
$ws.Range("D24").Value = 241
$ws.Range("E24").Value = 2.489626556016
$ws.Range("F24").Value = 930
$ws.Range("G24").Value = 895
$ws.Range("H24").Value = 3.910614525139
$ws.Range("I24").Value = 2741
$ws.Range("J24").Value = 3072
$ws.Range("K24").Value = -10.774739583333
$ws.Range("L24").Value = -17.786442711457
$ws.Range("M24").Value = 49.945295404814

# Row 25
$ws.Range("D25").Value = 139
$ws.Range("E25").Value = -7.194244604316
$ws.Range("F25").Value = 456
$ws.Range("G25").Value = 467
$ws.Range("H25").Value = -2.355460385438
$ws.Range("I25").Value = 1436
$ws.Range("J25").Value = 1741
$ws.Range("K25").Value = -17.518667432510
$ws.Range("L25").Value = -33.672055427251

# Row 26
$ws.Range("C26").Value = 86
$ws.Range("D26").Value = 104
$ws.Range("E26").Value = -17.307692307692
$ws.Range("F26").Value = 381
$ws.Range("G26").Value = 346
$ws.Range("H26").Value = 10.115606936416
$ws.Range("I26").Value = 1036
$ws.Range("J26").Value = 966
$ws.Range("K26").Value = 7.246376811594
$ws.Range("L26").Value = 5.284552845528
$ws.Range("M26").Value = -6.834532374100

# Row 27
$ws.Range("C27").Value = 1
$ws.Range("E27").Value = -66.666666666666
$ws.Range("F27").Value = 11
$ws.Range("G27").Value = 16
$ws.Range("H27").Value = -31.25
$ws.Range("I27").Value = 44
$ws.Range("J27").Value = 58
$ws.Range("K27").Value = -24.137931034482
$ws.Range("L27").Value = -24.137931034482

# Row 28
$ws.Range("C28").Value = 15
$ws.Range("D28").Value = 12
$ws.Range("E28").Value = 25
$ws.Range("F28").Value = 53
$ws.Range("G28").Value = 43
$ws.Range("H28").Value = 23.255813953488
$ws.Range("I28").Value = 129
$ws.Range("J28").Value = 127
$ws.Range("K28").Value = 1.574803149606
$ws.Range("L28").Value = -12.244897959183

# Row 29
$ws.Range("C29").Value = 3
$ws.Range("D29").Value = 1
$ws.Range("E29").Value = 200
$ws.Range("F29").Value = 9
$ws.Range("G29").Value = 11
$ws.Range("H29").Value = -18.181818181818
$ws.Range("I29").Value = 23
$ws.Range("J29").Value = 34
$ws.Range("K29").Value = -32.352941176470
$ws.Range("L29").Value = -42.5
$ws.Range("M29").Value = -34.285714285714
$ws.Range("N29").Value = -87.150837988826

# Row 30
$ws.Range("C30").Value = 2
$ws.Range("D30").Value = 1
$ws.Range("E30").Value = 100
$ws.Range("F30").Value = 7
$ws.Range("G30").Value = 10
$ws.Range("H30").Value = -30
$ws.Range("I30").Value = 18
$ws.Range("J30").Value = 32
$ws.Range("K30").Value = -43.75
$ws.Range("L30").Value = -50
$ws.Range("M30").Value = -45.454545454545
$ws.Range("N30").Value = -89.156626506024

# Row 31
$ws.Range("D31").Value = 1
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 12
$ws.Range("H31").Value = 20
$ws.Range("I31").Value = 22
$ws.Range("J31").Value = 18
$ws.Range("K31").Value = 22.222222222222
$ws.Range("L31").Value = -37.142857142857

# Row 33
$ws.Range("F33").Value = 3
$ws.Range("I33").Value = 3
$ws.Range("J33").Value = 6
$ws.Range("K33").Value = -50
$ws.Range("L33").Value = -62.5

# --- Cells that were text placeholders ("0" / "***.*") now carry real
#     numbers; restore the normal count / percent number formats so the
#     style id matches the other numeric cells in the row. ---
$ws.Range("C31").Value = 1
$ws.Range('C31').NumberFormat = '#,##0'

$ws.Range("C33").Value = 2
$ws.Range('C33').NumberFormat = '#,##0'
$ws.Range("D33").Value = 1
$ws.Range('D33').NumberFormat = '#,##0'
$ws.Range("E33").Value = 100
$ws.Range('E33').NumberFormat = '#,##0.0;"-"#,##0.0'
